$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.77
$ws.Range("H2").Value = 4.6
$ws.Range("O2").Value = 1.16
$ws.Range("P2").Value = 2.84
$ws.Range("Q2").Value = 1.48
$ws.Range("R2").Value = 1.74
$ws.Range("S2").Value = 2.18
$ws.Range("U2").Value = 2.64
$ws.Range("W2").Value = 2.3
$ws.Range("AG2").Value = 11
$ws.Range("AN2").Value = 7.6

# Row 3
$ws.Range("G3").Value = 2.04
$ws.Range("J3").Value = 3.8
$ws.Range("L3").Value = 1.32
$ws.Range("O3").Value = 1.22
$ws.Range("Q3").Value = 1.67
$ws.Range("W3").Value = 1.96
$ws.Range("X3").Value = 24
$ws.Range("AA3").Value = 100
$ws.Range("AI3").Value = 50
$ws.Range("AN3").Value = 11

# Row 6
$ws.Range("N6").Value = 3.7
$ws.Range("O6").Value = 1.35
$ws.Range("Q6").Value = 2.04
$ws.Range("T6").Value = 1.83
$ws.Range("AE6").Value = 44

# Row 7
$ws.Range("N7").Value = 5.3
$ws.Range("S7").Value = 2.66

# Row 9
$ws.Range("G9").Value = 2.42
$ws.Range("O9").Value = 1.3
$ws.Range("X9").Value = 15.5
$ws.Range("AE9").Value = 36

# Row 10
$ws.Range("H10").Value = 1.47
$ws.Range("I10").Value = 1.48
